# Convert the legacy "{$ img:<key> $}" image placeholder syntax to the
# standard Jinja2 "{{ <key> }}" syntax used by docxtpl's InlineImage.

$d = $word.ActiveDocument

# {$ img:cad_model $}  ->  {{ cad_model }}
$d.Content.Find.Execute("{`$ img:cad_model `$}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ cad_model }}", 2)

# {$ img:plot $}  ->  {{ plot }}
$d.Content.Find.Execute("{`$ img:plot `$}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ plot }}", 2)
